$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1154.6364
$ws.Range("J40").Value = 1200.25
$ws.Range("L40").Value = 1200.25
$ws.Range("N40").Value = -1550.25
$ws.Range("H70").Value = 3374.25
$ws.Range("J70").Value = 3332.6667
$ws.Range("L70").Value = 9998.000100000001
$ws.Range("N70").Value = -10538.0001
$ws.Range("H73").Value = 3374.25
$ws.Range("J73").Value = 3332.6667
$ws.Range("L73").Value = 9998.000100000001
$ws.Range("N73").Value = -11870.0001
$ws.Range("H86").Value = 23333.334
$ws.Range("I86").Value = 20000
$ws.Range("K86").Value = 20000
$ws.Range("M86").Value = -18877
$ws.Range("H89").Value = 23333.334
$ws.Range("I89").Value = 20000
$ws.Range("K89").Value = 100000
$ws.Range("M89").Value = -94384
$ws.Range("H106").Value = 9686.857
$ws.Range("I106").Value = 9686.857
$ws.Range("K106").Value = 9686.857
$ws.Range("M106").Value = -9055.857
$ws.Range("H116").Value = 41218.625
$ws.Range("I116").Value = 24125
$ws.Range("J116").Value = 46916.5
$ws.Range("K116").Value = 24125
$ws.Range("L116").Value = 46916.5
$ws.Range("M116").Value = -20683
$ws.Range("N116").Value = -53800.5
$ws.Range("H132").Value = 51037.25
$ws.Range("J132").Value = 1382
$ws.Range("L132").Value = 4146
$ws.Range("N132").Value = -9206
$ws.Range("H135").Value = 8977.916999999999
$ws.Range("I135").Value = 8810
$ws.Range("K135").Value = 79290
$ws.Range("M135").Value = -76755
$ws.Range("H137").Value = 68333
$ws.Range("I137").Value = 99999
$ws.Range("K137").Value = 299997
$ws.Range("M137").Value = -297447
$ws.Range("H138").Value = 4302.485
$ws.Range("I138").Value = 7311.25
$ws.Range("J138").Value = 3339.68
$ws.Range("K138").Value = 21933.75
$ws.Range("L138").Value = 10019.04
$ws.Range("M138").Value = -16793.75
$ws.Range("N138").Value = -20299.04
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2416.25
$ws.Range("I2").Value = 2295.1428
$ws.Range("J2").Value = 2585.8
$ws.Range("K2").Value = 2295.1428
$ws.Range("L2").Value = 2585.8
$ws.Range("M2").Value = -2182.1428
$ws.Range("N2").Value = -2811.8
$ws.Range("H45").Value = 4220
$ws.Range("I45").Value = 3775
$ws.Range("K45").Value = 3775
$ws.Range("M45").Value = -3398
$ws.Range("H74").Value = 5365.3213
$ws.Range("I74").Value = 2538.8298
$ws.Range("K74").Value = 2538.8298
$ws.Range("M74").Value = -1664.8298
$ws.Range("H77").Value = 5365.3213
$ws.Range("I77").Value = 2538.8298
$ws.Range("K77").Value = 12694.149
$ws.Range("M77").Value = -8326.148999999999
$ws.Range("H116").Value = 2416.25
$ws.Range("I116").Value = 2295.1428
$ws.Range("J116").Value = 2585.8
$ws.Range("K116").Value = 2295.1428
$ws.Range("L116").Value = 2585.8
$ws.Range("M116").Value = -1.142800000000079
$ws.Range("N116").Value = -7173.8
$ws.Range("H122").Value = 71429540
$ws.Range("I122").Value = 100000430
$ws.Range("K122").Value = 300001290
$ws.Range("M122").Value = -299998840
$ws.Range("H133").Value = 71999.5
$ws.Range("J133").Value = 71999.5
$ws.Range("L133").Value = 71999.5
$ws.Range("N133").Value = -77059.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2416.25
$ws.Range("I3").Value = 2295.1428
$ws.Range("J3").Value = 2585.8
$ws.Range("K3").Value = 2295.1428
$ws.Range("L3").Value = 2585.8
$ws.Range("M3").Value = -2181.1428
$ws.Range("N3").Value = -2813.8
$ws.Range("H20").Value = 48828.363
$ws.Range("I20").Value = 80847.38
$ws.Range("K20").Value = 80847.38
$ws.Range("M20").Value = -80600.38
$ws.Range("H82").Value = 30745.75
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766
$ws.Range("H85").Value = 30745.75
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2941906
$ws.Range("I19").Value = 3846907.8
$ws.Range("J19").Value = 650
$ws.Range("K19").Value = 3846907.8
$ws.Range("L19").Value = 650
$ws.Range("M19").Value = -3846737.8
$ws.Range("N19").Value = -990
$ws.Range("H24").Value = 2941906
$ws.Range("I24").Value = 3846907.8
$ws.Range("J24").Value = 650
$ws.Range("K24").Value = 3846907.8
$ws.Range("L24").Value = 650
$ws.Range("M24").Value = -3846737.8
$ws.Range("N24").Value = -990
$ws.Range("H31").Value = 2274.158
$ws.Range("I31").Value = 3243.3333
$ws.Range("J31").Value = 1401.9
$ws.Range("K31").Value = 3243.3333
$ws.Range("L31").Value = 1401.9
$ws.Range("M31").Value = -2948.3333
$ws.Range("N31").Value = -1991.9
$ws.Range("H34").Value = 2274.158
$ws.Range("I34").Value = 3243.3333
$ws.Range("J34").Value = 1401.9
$ws.Range("K34").Value = 3243.3333
$ws.Range("L34").Value = 1401.9
$ws.Range("M34").Value = -3041.3333
$ws.Range("N34").Value = -1805.9
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -18740
$ws.Range("H68").Value = 31436
$ws.Range("J68").Value = 31436
$ws.Range("L68").Value = 31436
$ws.Range("N68").Value = -32934
$ws.Range("H71").Value = 31436
$ws.Range("J71").Value = 31436
$ws.Range("L71").Value = 94308
$ws.Range("N71").Value = -101796
$ws.Range("H132").Value = 3395.7778
$ws.Range("I132").Value = 3395.7778
$ws.Range("K132").Value = 10187.3334
$ws.Range("M132").Value = -7657.3334
$ws.Range("H134").Value = 2344.9524
$ws.Range("I134").Value = 2344.9524
$ws.Range("K134").Value = 7034.8572
$ws.Range("M134").Value = -4499.8572
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 531.2857
$ws.Range("I86").Value = 556.75
$ws.Range("J86").Value = 497.33334
$ws.Range("K86").Value = 1670.25
$ws.Range("L86").Value = 1492.00002
$ws.Range("M86").Value = -484.25
$ws.Range("N86").Value = -3864.00002
$ws.Range("H89").Value = 531.2857
$ws.Range("I89").Value = 556.75
$ws.Range("J89").Value = 497.33334
$ws.Range("K89").Value = 5010.75
$ws.Range("L89").Value = 4476.00006
$ws.Range("M89").Value = 917.25
$ws.Range("N89").Value = -16332.00006
$ws.Range("H113").Value = 801.04346
$ws.Range("I113").Value = 532.6667
$ws.Range("J113").Value = 841.3
$ws.Range("K113").Value = 1598.0001
$ws.Range("L113").Value = 2523.9
$ws.Range("M113").Value = 571.9999
$ws.Range("N113").Value = -6863.9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3206.1667
$ws.Range("I80").Value = 3047.4
$ws.Range("K80").Value = 3047.4
$ws.Range("M80").Value = -2049.4
$ws.Range("H83").Value = 3206.1667
$ws.Range("I83").Value = 3047.4
$ws.Range("K83").Value = 15237
$ws.Range("M83").Value = -10245
$ws.Range("H122").Value = 3362.9583
$ws.Range("I122").Value = 3277.1177
$ws.Range("K122").Value = 9831.3531
$ws.Range("M122").Value = -7381.3531
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1087.7826
$ws.Range("J55").Value = 1234.8
$ws.Range("L55").Value = 1234.8
$ws.Range("N55").Value = -1580.8
$ws.Range("H61").Value = 6825.2793
$ws.Range("I61").Value = 6622.4863
$ws.Range("J61").Value = 8075.8335
$ws.Range("K61").Value = 6622.4863
$ws.Range("L61").Value = 8075.8335
$ws.Range("M61").Value = -6420.4863
$ws.Range("N61").Value = -8479.833500000001
$ws.Range("H68").Value = 8144.3335
$ws.Range("I68").Value = 6912.5
$ws.Range("J68").Value = 17999
$ws.Range("K68").Value = 6912.5
$ws.Range("L68").Value = 17999
$ws.Range("M68").Value = -6163.5
$ws.Range("N68").Value = -19497
$ws.Range("H71").Value = 8144.3335
$ws.Range("I71").Value = 6912.5
$ws.Range("J71").Value = 17999
$ws.Range("K71").Value = 34562.5
$ws.Range("L71").Value = 89995
$ws.Range("M71").Value = -30818.5
$ws.Range("N71").Value = -97483
$ws.Range("H103").Value = 25499.5
$ws.Range("J103").Value = 25499.5
$ws.Range("L103").Value = 25499.5
$ws.Range("N103").Value = -27843.5
$ws.Range("H113").Value = 6825.2793
$ws.Range("I113").Value = 6622.4863
$ws.Range("J113").Value = 8075.8335
$ws.Range("K113").Value = 6622.4863
$ws.Range("L113").Value = 8075.8335
$ws.Range("M113").Value = -4452.4863
$ws.Range("N113").Value = -12415.8335
$ws.Range("H122").Value = 4837.8667
$ws.Range("I122").Value = 3312.9
$ws.Range("J122").Value = 7887.8
$ws.Range("K122").Value = 9938.700000000001
$ws.Range("L122").Value = 23663.4
$ws.Range("M122").Value = -7488.700000000001
$ws.Range("N122").Value = -28563.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 20013918
$ws.Range("I74").Value = 100000000
$ws.Range("J74").Value = 17397.25
$ws.Range("K74").Value = 100000000
$ws.Range("L74").Value = 17397.25
$ws.Range("M74").Value = -99999064
$ws.Range("N74").Value = -19269.25
$ws.Range("H77").Value = 20013918
$ws.Range("I77").Value = 100000000
$ws.Range("J77").Value = 17397.25
$ws.Range("K77").Value = 300000000
$ws.Range("L77").Value = 52191.75
$ws.Range("M77").Value = -61551.75
$ws.Range("N77").Value = -61551.75
$ws.Range("H122").Value = 44722.703
$ws.Range("I122").Value = 3099.8635
$ws.Range("K122").Value = 9299.5905
$ws.Range("M122").Value = -6849.5905
$ws.Range("H126").Value = 1499.8695
$ws.Range("I126").Value = 1399.65
$ws.Range("J126").Value = 2168
$ws.Range("K126").Value = 4198.950000000001
$ws.Range("L126").Value = 6504
$ws.Range("M126").Value = -1728.950000000001
$ws.Range("N126").Value = -11444
